$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet's conversion note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.09 = 37783.54 pesos`n✅ 37783.54 pesos = 9.0 = 949.54 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet's rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 109.995
$ws2.Range("O10").Value = 4156
$ws2.Range("N12").Value = 4198
$ws2.Range("O12").Value = 105.5
